# Commit: "fix error when server exits"
# On the "Problems" slide (sldId=455), update the bullet list text:
#   - "The client apps do not exit automatically when the server shuts down "
#     becomes "The client apps can't handle special characters "
#   - "The maximum number of clients has not been quantified yet"
#     becomes "The maximum number of clients has not been quantified"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para2 = $tr.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "The client apps can't handle special characters "

$para3 = $tr.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = "The maximum number of clients has not been quantified"
